# Fix the "2050" column header label on sheets 1-5 (it was mistakenly left
# as a stray numeric value instead of the text label used by the other
# header cells), and drop the "Total" summary row that was appended to the
# bottom of every table.

$wb = $excel.ActiveWorkbook

# Sheet 1: "Potencia Acumulada - SIN (MW)"
$ws = $wb.Worksheets.Item(1)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"
$ws.Rows.Item(13).Delete()

# Sheet 2: "Geracao Periodo Medio (MWMed)"
$ws = $wb.Worksheets.Item(2)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"
$ws.Rows.Item(13).Delete()

# Sheet 3: "Atendimento a Ponta(MW)"
$ws = $wb.Worksheets.Item(3)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"
$ws.Rows.Item(13).Delete()

# Sheet 4: "Potencia Incremental - SIN(MW)" uses year-range labels, so the
# corrected header is "2041-2050" instead of a plain year.
$ws = $wb.Worksheets.Item(4)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2041-2050"
$ws.Rows.Item(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" only needs the header fix; it never
# had a Total row.
$ws = $wb.Worksheets.Item(5)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"

# Sheet 6: "Custo Total (bilhoes de R$)" has its own Total row at row 4
# (no header label issue on this sheet).
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(4).Delete()
